# "Generate Report for Archive" — re-sort the 4 tracked-file rows (rows 2-5)
# on every sheet into the new report order:
#   41ac75ec-...  ->  d17ea57b-...  ->  dc34049b-...  ->  12d4c890-...
# The first three files are now "In Translation"; the last keeps
# "Ready for handoff". Hyperlinks for the re-ordered rows are rebuilt so the
# displayed file name always matches the row it sits in, while the
# hyperlink's target URL stays pinned to the same underlying file it always
# pointed to.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "41ac75ec-abd2-4957-b215-db8c313ea693.md"
$ws1.Range("B2").Value = "In Translation"
$ws1.Range("C2").Value = "In Translation"

$ws1.Range("A3").Value = "d17ea57b-83ce-450a-b525-29e57a9e9596.md"
$ws1.Range("B3").Value = "In Translation"
$ws1.Range("C3").Value = "In Translation"

$ws1.Range("A4").Value = "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"

$ws1.Range("A5").Value = "12d4c890-b096-457c-8a96-b5e95ef57d99.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe83648533a2d131c1552e15f3094ed31fdd680e/e2e/41ac75ec-abd2-4957-b215-db8c313ea693.md", "", "", "41ac75ec-abd2-4957-b215-db8c313ea693.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/d17ea57b-83ce-450a-b525-29e57a9e9596.md", "", "", "d17ea57b-83ce-450a-b525-29e57a9e9596.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/dc34049b-f3d8-43c1-babc-2a17df54f6e2.md", "", "", "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/e2e/12d4c890-b096-457c-8a96-b5e95ef57d99.md", "", "", "12d4c890-b096-457c-8a96-b5e95ef57d99.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "41ac75ec-abd2-4957-b215-db8c313ea693.md"
$ws2.Range("B2").Value = "In Translation"
$ws2.Range("C2").Value = "41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-06 03:37:49"

$ws2.Range("A3").Value = "d17ea57b-83ce-450a-b525-29e57a9e9596.md"
$ws2.Range("B3").Value = "In Translation"
$ws2.Range("C3").Value = "d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-06 03:35:59"

$ws2.Range("A4").Value = "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md"
$ws2.Range("B4").Value = "In Translation"
$ws2.Range("C4").Value = "dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-02-06 03:35:59"

$ws2.Range("A5").Value = "12d4c890-b096-457c-8a96-b5e95ef57d99.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-02-06 03:38:30"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe83648533a2d131c1552e15f3094ed31fdd680e/e2e/41ac75ec-abd2-4957-b215-db8c313ea693.md", "", "", "41ac75ec-abd2-4957-b215-db8c313ea693.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca0e89ba70701e04fa465f7bda0d630947ba123b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.zh-cn.xlf", "", "", "41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/d17ea57b-83ce-450a-b525-29e57a9e9596.md", "", "", "d17ea57b-83ce-450a-b525-29e57a9e9596.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f349ccbcfea23be95b4f01f5c47bc8ef5c941e79/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.zh-cn.xlf", "", "", "d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/dc34049b-f3d8-43c1-babc-2a17df54f6e2.md", "", "", "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f349ccbcfea23be95b4f01f5c47bc8ef5c941e79/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.zh-cn.xlf", "", "", "dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/e2e/12d4c890-b096-457c-8a96-b5e95ef57d99.md", "", "", "12d4c890-b096-457c-8a96-b5e95ef57d99.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2640b39a149cdb555e6779efa20f44bf9554f4d2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.zh-cn.xlf", "", "", "12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "41ac75ec-abd2-4957-b215-db8c313ea693.md"
$ws3.Range("B2").Value = "In Translation"
$ws3.Range("C2").Value = "41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-06 03:38:00"

$ws3.Range("A3").Value = "d17ea57b-83ce-450a-b525-29e57a9e9596.md"
$ws3.Range("B3").Value = "In Translation"
$ws3.Range("C3").Value = "d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-06 03:36:28"

$ws3.Range("A4").Value = "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md"
$ws3.Range("B4").Value = "In Translation"
$ws3.Range("C4").Value = "dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.de-de.xlf"
$ws3.Range("D4").Value = "2016-02-06 03:36:28"

$ws3.Range("A5").Value = "12d4c890-b096-457c-8a96-b5e95ef57d99.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.de-de.xlf"
$ws3.Range("D5").Value = "2016-02-06 03:38:41"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fe83648533a2d131c1552e15f3094ed31fdd680e/e2e/41ac75ec-abd2-4957-b215-db8c313ea693.md", "", "", "41ac75ec-abd2-4957-b215-db8c313ea693.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/764a1972a5ca595d0ee08303b6ee0baddb76cdcb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.de-de.xlf", "", "", "41ac75ec-abd2-4957-b215-db8c313ea693.10d08a91af024f3222b76b7176bb2f6d6e4482ee.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/d17ea57b-83ce-450a-b525-29e57a9e9596.md", "", "", "d17ea57b-83ce-450a-b525-29e57a9e9596.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7fe4738ec0cda41d4ff0bb775d0df6e9b15aaa0c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.de-de.xlf", "", "", "d17ea57b-83ce-450a-b525-29e57a9e9596.c4de3d7430aba81ed21bbae71e02c6de232632ad.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11a9d22124ae7249383a3741c1ea28827df430cb/e2e/dc34049b-f3d8-43c1-babc-2a17df54f6e2.md", "", "", "dc34049b-f3d8-43c1-babc-2a17df54f6e2.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7fe4738ec0cda41d4ff0bb775d0df6e9b15aaa0c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.de-de.xlf", "", "", "dc34049b-f3d8-43c1-babc-2a17df54f6e2.2b9b782ce1559c126678302a5892deec8bd13720.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/e2e/12d4c890-b096-457c-8a96-b5e95ef57d99.md", "", "", "12d4c890-b096-457c-8a96-b5e95ef57d99.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/016b2a6ad753749ab9e501cd37a38d7be4a5e3b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.de-de.xlf", "", "", "12d4c890-b096-457c-8a96-b5e95ef57d99.897a9afa869adaa1eee29b7bdec54343a9fbce9a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/6f6483c0a2825786867df2f6c4f640503bf98c15/.localization-config", "", "", ".localization-config")
